# Update gh-pages output (commit 456a3b4)
# Applies cell-value refresh to 展览(1), 演出(2), 全部类型(4) sheets:
#  - F column ("想去人数" / want-to-go count) values refreshed from live site
#  - I column ("Cover" image URL) refreshed for rows whose cover image changed

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# --- 展览 sheet1 (Worksheets.Item(1)) ---
$ws1.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202407/h8BarH1S1720583334377.jpeg"
$ws1.Range("F4").Value = 0
$ws1.Range("F5").Value = 1162
$ws1.Range("F11").Value = 0
$ws1.Range("F12").Value = 0
$ws1.Range("F13").Value = 0
$ws1.Range("F14").Value = 0
$ws1.Range("F15").Value = 0
$ws1.Range("F16").Value = 0
$ws1.Range("F17").Value = 412
$ws1.Range("F18").Value = 0
$ws1.Range("F19").Value = 586
$ws1.Range("F21").Value = 211
$ws1.Range("F22").Value = 161
$ws1.Range("F23").Value = 10194
$ws1.Range("F27").Value = 0
$ws1.Range("F31").Value = 0
$ws1.Range("F33").Value = 17
$ws1.Range("F34").Value = 0
$ws1.Range("F37").Value = 1406
$ws1.Range("F39").Value = 0
$ws1.Range("F44").Value = 1103
$ws1.Range("F46").Value = 0
$ws1.Range("F47").Value = 0
$ws1.Range("F48").Value = 67

# --- 演出 sheet2 (Worksheets.Item(2)) ---
$ws2.Range("F2").Value = 0
$ws2.Range("F3").Value = 0
$ws2.Range("F8").Value = 7
$ws2.Range("F9").Value = 42
$ws2.Range("F11").Value = 1
$ws2.Range("F14").Value = 5
$ws2.Range("F15").Value = 97
$ws2.Range("F17").Value = 0

# --- 全部类型 sheet4 (Worksheets.Item(4)) ---
$ws4.Range("F2").Value = 40
$ws4.Range("F3").Value = 36
$ws4.Range("F4").Value = 72
$ws4.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202407/h8BarH1S1720583334377.jpeg"
$ws4.Range("F5").Value = 27
$ws4.Range("F8").Value = 9056
$ws4.Range("F10").Value = 247
$ws4.Range("F11").Value = 7159
$ws4.Range("F15").Value = 3
$ws4.Range("F16").Value = 5522
$ws4.Range("F17").Value = 73
$ws4.Range("F18").Value = 6308
$ws4.Range("F19").Value = 6308
$ws4.Range("F21").Value = 423
$ws4.Range("F24").Value = 277
$ws4.Range("F25").Value = 211
$ws4.Range("F26").Value = 0
$ws4.Range("F27").Value = 10194
$ws4.Range("F32").Value = 84
$ws4.Range("F37").Value = 2101
$ws4.Range("F38").Value = 312
$ws4.Range("F39").Value = 1406
$ws4.Range("F42").Value = 0
$ws4.Range("F43").Value = 0
$ws4.Range("F45").Value = 1103
$ws4.Range("F46").Value = 1078
$ws4.Range("F49").Value = 67
$ws4.Range("F50").Value = 1096
